$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet1: "custom-formula" -> "custom function"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new banner row at the very top (everything else shifts down
# by one row - formulas/refs adjust automatically).
$ws1.Rows.Item(1).Insert() | Out-Null

# New row 1: a big merged "Custom Function" title.
$ws1.Range("A1").Value = "Custom Function"
$ws1.Range("A1:C1").Merge() | Out-Null
$ws1.Range("A1").Font.Size = 22
$ws1.Range("A1").HorizontalAlignment = -4108
$ws1.Rows.Item(1).RowHeight = 29

# Row 2 (previously row 1): re-label two of the header cells and bump
# the header font up to match the new title's emphasis.
$ws1.Range("B2").Value = "Custom function"
$ws1.Range("A2").Value = "Function Name"
$ws1.Range("A2:G2").Font.Size = 14

# Sheet is now renamed last, once its own edits are in place (the
# rename forces a workbook-wide formula refresh).
$ws1.Name = "custom function"

# ------------------------------------------------------------------
# Active sheet / selections
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B8").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B4").Select() | Out-Null
